# hw5/Statistics.xlsx - "Updated the source code for problem 2 to use what I
# think is the correct trace options to get the basic block count."
#
# The only sheet with data changes is "Problem2" (re-run trace numbers for
# the basic-block count). "Problem1" is untouched content-wise.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Problem2")

# --- Row 2 ("621_wrf_s train"): Basic Blocks + Total Instructions updated ---
$ws.Range("B2").Value = 3104
$ws.Range("E2").Value = 482639296

# --- Row 3 ("623_xalanchbmk_s train"): Basic Blocks + Total Instructions ---
$ws.Range("B3").Value = 4619
$ws.Range("E3").Value = 1421192

# --- Row 4 ("mm_mult_serial 256x256"): Basic Blocks, Memory Reads/Writes,
#     Total Instructions all updated ---
$ws.Range("B4").Value = 104862
$ws.Range("C4").Value = 260269
$ws.Range("D4").Value = 251753
$ws.Range("E4").Value = 583459817334

# --- Row 5 ("accumulate array 1x256"): Basic Blocks, Memory Reads, Total
#     Instructions updated; this row also picks up a distinct font (it's the
#     last-touched row, highlighted with its own font record in styles.xml) ---
$ws.Range("B5").Value = 37987
$ws.Range("C5").Value = 53881
$ws.Range("E5").Value = 257156226729

# Give row 5 its own font (keeps the same border/number-format, just a
# fresh font record) to mirror the new cellXfs entries introduced upstream.
$ws.Range("A5:E5").Font.Name = "Calibri"

# The author's last selection before saving ended up on E5 of Problem2.
$ws.Range("E5").Select() | Out-Null
